$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2033.1666
$ws.Cells.Item(40, 9).Value = 1939.8
$ws.Cells.Item(40, 10).Value = 2500
$ws.Cells.Item(40, 11).Value = 1939.8
$ws.Cells.Item(40, 12).Value = 2500
$ws.Cells.Item(40, 13).Value = -1764.8
$ws.Cells.Item(40, 14).Value = -2850

# ALC row 47
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(47, 8).Value = 70000
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 70000
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 70000
$ws.Cells.Item(47, 13).ClearContents()
$ws.Cells.Item(47, 14).Value = -71944

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).ClearContents()

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 632.9474
$ws.Cells.Item(107, 9).Value = 648.64703
$ws.Cells.Item(107, 10).Value = 499.5
$ws.Cells.Item(107, 11).Value = 648.64703
$ws.Cells.Item(107, 12).Value = 499.5
$ws.Cells.Item(107, 13).Value = 1271.35297
$ws.Cells.Item(107, 14).Value = -4339.5

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4285.091
$ws.Cells.Item(116, 9).Value = 4289.2
$ws.Cells.Item(116, 10).Value = 4281.6665
$ws.Cells.Item(116, 11).Value = 4289.2
$ws.Cells.Item(116, 12).Value = 4281.6665
$ws.Cells.Item(116, 13).Value = -847.1999999999998
$ws.Cells.Item(116, 14).Value = -11165.6665

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1110.1666
$ws.Cells.Item(137, 9).Value = 890.13336
$ws.Cells.Item(137, 10).Value = 2210.3333
$ws.Cells.Item(137, 11).Value = 2670.40008
$ws.Cells.Item(137, 12).Value = 6630.999899999999
$ws.Cells.Item(137, 13).Value = -120.4000800000003
$ws.Cells.Item(137, 14).Value = -11730.9999

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 721.9
$ws.Cells.Item(2, 9).Value = 721.9
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 721.9
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -608.9

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 474.16666
$ws.Cells.Item(4, 9).Value = 474.16666
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 474.16666
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -358.16666
$ws.Cells.Item(4, 14).ClearContents()

# ARM row 35
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 2075
$ws.Cells.Item(35, 9).Value = 2075
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 2075
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -1669

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1845.6666
$ws.Cells.Item(74, 9).Value = 1845.6666
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 1845.6666
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -971.6666

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1845.6666
$ws.Cells.Item(77, 9).Value = 1845.6666
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 9228.333000000001
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -4860.333000000001

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 60528.824
$ws.Cells.Item(102, 9).Value = 72785
$ws.Cells.Item(102, 10).Value = 3333.3333
$ws.Cells.Item(102, 11).Value = 72785
$ws.Cells.Item(102, 12).Value = 3333.3333
$ws.Cells.Item(102, 13).Value = -71163
$ws.Cells.Item(102, 14).Value = -6577.3333

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 795
$ws.Cells.Item(110, 9).Value = 795
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 795
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = 1250

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 721.9
$ws.Cells.Item(116, 9).Value = 721.9
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 721.9
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 1572.1

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1849
$ws.Cells.Item(122, 9).Value = 1832.3334
$ws.Cells.Item(122, 10).Value = 1899
$ws.Cells.Item(122, 11).Value = 5497.0002
$ws.Cells.Item(122, 12).Value = 5697
$ws.Cells.Item(122, 13).Value = -3047.0002
$ws.Cells.Item(122, 14).Value = -10597

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 721.9
$ws.Cells.Item(3, 9).Value = 721.9
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 721.9
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -607.9

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1014.2
$ws.Cells.Item(20, 9).Value = 819.8570999999999
$ws.Cells.Item(20, 10).Value = 1467.6666
$ws.Cells.Item(20, 11).Value = 819.8570999999999
$ws.Cells.Item(20, 12).Value = 1467.6666
$ws.Cells.Item(20, 13).Value = -572.8570999999999
$ws.Cells.Item(20, 14).Value = -1961.6666

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1058.6
$ws.Cells.Item(94, 9).Value = 916.9474
$ws.Cells.Item(94, 10).Value = 3750
$ws.Cells.Item(94, 11).Value = 916.9474
$ws.Cells.Item(94, 12).Value = 3750
$ws.Cells.Item(94, 13).Value = -465.9474
$ws.Cells.Item(94, 14).Value = -4652

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3809.1155
$ws.Cells.Item(99, 9).Value = 3523.4348
$ws.Cells.Item(99, 10).Value = 5999.3335
$ws.Cells.Item(99, 11).Value = 3523.4348
$ws.Cells.Item(99, 12).Value = 5999.3335
$ws.Cells.Item(99, 13).Value = -2025.4348
$ws.Cells.Item(99, 14).Value = -8995.333500000001

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1563.3334
$ws.Cells.Item(105, 9).Value = 1563.3334
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 1563.3334
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = 183.6666

# CRP row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 829.3333
$ws.Cells.Item(10, 9).Value = 829.3333
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 829.3333
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = -690.3333
$ws.Cells.Item(10, 14).ClearContents()

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2982.35
$ws.Cells.Item(31, 9).Value = 1617.2667
$ws.Cells.Item(31, 10).Value = 7077.6
$ws.Cells.Item(31, 11).Value = 1617.2667
$ws.Cells.Item(31, 12).Value = 7077.6
$ws.Cells.Item(31, 13).Value = -1322.2667
$ws.Cells.Item(31, 14).Value = -7667.6

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2982.35
$ws.Cells.Item(34, 9).Value = 1617.2667
$ws.Cells.Item(34, 10).Value = 7077.6
$ws.Cells.Item(34, 11).Value = 1617.2667
$ws.Cells.Item(34, 12).Value = 7077.6
$ws.Cells.Item(34, 13).Value = -1415.2667
$ws.Cells.Item(34, 14).Value = -7481.6

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2036.25
$ws.Cells.Item(58, 9).Value = 2048.3333
$ws.Cells.Item(58, 10).Value = 2000
$ws.Cells.Item(58, 11).Value = 2048.3333
$ws.Cells.Item(58, 12).Value = 2000
$ws.Cells.Item(58, 13).Value = -1845.3333
$ws.Cells.Item(58, 14).Value = -2406

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2036.25
$ws.Cells.Item(136, 9).Value = 2048.3333
$ws.Cells.Item(136, 10).Value = 2000
$ws.Cells.Item(136, 11).Value = 6144.999899999999
$ws.Cells.Item(136, 12).Value = 6000
$ws.Cells.Item(136, 13).Value = -3594.999899999999
$ws.Cells.Item(136, 14).Value = -11100

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1323.4445
$ws.Cells.Item(4, 9).Value = 1264.9565
$ws.Cells.Item(4, 10).Value = 1659.75
$ws.Cells.Item(4, 11).Value = 3794.8695
$ws.Cells.Item(4, 12).Value = 4979.25
$ws.Cells.Item(4, 13).Value = -3682.8695
$ws.Cells.Item(4, 14).Value = -5203.25

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 267.22223
$ws.Cells.Item(12, 9).Value = 2.75
$ws.Cells.Item(12, 10).Value = 478.8
$ws.Cells.Item(12, 11).Value = 8.25
$ws.Cells.Item(12, 12).Value = 1436.4
$ws.Cells.Item(12, 13).Value = 164.75
$ws.Cells.Item(12, 14).Value = -1782.4

# CUL row 35
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(35, 8).Value = 500
$ws.Cells.Item(35, 9).Value = 500
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 1500
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -1212
$ws.Cells.Item(35, 14).ClearContents()

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 62752
$ws.Cells.Item(5, 9).Value = 62752
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 62752
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -62640

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 799.1429000000001
$ws.Cells.Item(97, 9).Value = 333
$ws.Cells.Item(97, 10).Value = 1148.75
$ws.Cells.Item(97, 11).Value = 333
$ws.Cells.Item(97, 12).Value = 1148.75
$ws.Cells.Item(97, 13).Value = 163
$ws.Cells.Item(97, 14).Value = -2140.75

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 3875
$ws.Cells.Item(2, 9).Value = 500
$ws.Cells.Item(2, 10).Value = 5000
$ws.Cells.Item(2, 11).Value = 500
$ws.Cells.Item(2, 12).Value = 5000
$ws.Cells.Item(2, 13).Value = -388
$ws.Cells.Item(2, 14).Value = -5224

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 144.38461
$ws.Cells.Item(46, 9).Value = 143.18182
$ws.Cells.Item(46, 10).Value = 151
$ws.Cells.Item(46, 11).Value = 143.18182
$ws.Cells.Item(46, 12).Value = 151
$ws.Cells.Item(46, 13).Value = 44.81818000000001
$ws.Cells.Item(46, 14).Value = -527

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 5064.3335
$ws.Cells.Item(61, 9).Value = 4511.4287
$ws.Cells.Item(61, 10).Value = 6999.5
$ws.Cells.Item(61, 11).Value = 4511.4287
$ws.Cells.Item(61, 12).Value = 6999.5
$ws.Cells.Item(61, 13).Value = -4309.4287
$ws.Cells.Item(61, 14).Value = -7403.5

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 5064.3335
$ws.Cells.Item(113, 9).Value = 4511.4287
$ws.Cells.Item(113, 10).Value = 6999.5
$ws.Cells.Item(113, 11).Value = 4511.4287
$ws.Cells.Item(113, 12).Value = 6999.5
$ws.Cells.Item(113, 13).Value = -2341.4287
$ws.Cells.Item(113, 14).Value = -11339.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3437.7222
$ws.Cells.Item(122, 9).Value = 2324.8572
$ws.Cells.Item(122, 10).Value = 4145.909
$ws.Cells.Item(122, 11).Value = 6974.571599999999
$ws.Cells.Item(122, 12).Value = 12437.727
$ws.Cells.Item(122, 13).Value = -4524.571599999999
$ws.Cells.Item(122, 14).Value = -17337.727

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 105
$ws.Cells.Item(2, 9).Value = 105
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 105
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 7
$ws.Cells.Item(2, 14).ClearContents()

# WVR row 28
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 14).ClearContents()
